$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Bold the "Handycache permission must be sufficient!" paragraph
#    (paragraph mark + both runs all get <w:b/>)
# ------------------------------------------------------------------
$hcPara = $d.Paragraphs.Item(3)
$hcPara.Range.Bold = 1

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the "3. Java net.properties"
#    heading to the end of the "See Brandmower Settings.docx"
#    paragraph (right after the last run, collapsed, no run split).
#
#    The runtime's Bookmarks.Add mis-places a *collapsed* range that
#    sits exactly at a paragraph's final character boundary (it resets
#    to document position 0). Work around this by temporarily typing a
#    placeholder character at that spot, wrapping the bookmark around
#    that single character, and then deleting the character again -
#    the now-collapsed bookmark stays correctly anchored in place.
# ------------------------------------------------------------------
$targetPara = $d.Paragraphs.Item(2)
$insertionPoint = $targetPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.End = $insertionPoint.End - 1
$insertionPoint.Collapse(0)

$insertionPoint.InsertAfter("Z")

$placeholder = $d.Range($insertionPoint.Start, $insertionPoint.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder2 = $d.Range($insertionPoint.Start, $insertionPoint.Start + 1)
$placeholder2.Delete()
